$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.720897
$ws.Range("H2").Value = 2.162691
$ws.Range("I2").Value = 0.0284720950782092
$ws.Range("J2").Value = 0.02847209507820921
$ws.Range("M2").Value = 0.8908616666666668
$ws.Range("N2").Value = 2.672585
$ws.Range("O2").Value = 0.04079002072021364
$ws.Range("P2").Value = 0.04079002072021363
$ws.Range("Q2").Value = 0.6422195029150001
$ws.Range("R2").Value = 5.779975526235001
$ws.Range("S2").Value = 0.001161377348188046
$ws.Range("T2").Value = 0.001161377348188046
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.720897
$ws.Range("H3").Value = 2.162691
$ws.Range("I3").Value = 0.0284720950782092
$ws.Range("J3").Value = 0.02847209507820921
$ws.Range("M3").Value = 16.81477433333333
$ws.Range("O3").Value = 0.7699006693471485
$ws.Range("P3").Value = 0.7699006693471484
$ws.Range("Q3").Value = 12.121720372577
$ws.Range("R3").Value = 109.095483353193
$ws.Range("S3").Value = 0.02192068505842892
$ws.Range("T3").Value = 0.02192068505842892
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.720897
$ws.Range("H4").Value = 2.162691
$ws.Range("I4").Value = 0.0284720950782092
$ws.Range("J4").Value = 0.02847209507820921
$ws.Range("M4").Value = 3.879966
$ws.Range("N4").Value = 11.639898
$ws.Range("O4").Value = 0.1776526024808091
$ws.Range("P4").Value = 0.1776526024808091
$ws.Range("Q4").Value = 2.797055849502
$ws.Range("R4").Value = 25.173502645518
$ws.Range("S4").Value = 0.005058141788724901
$ws.Range("T4").Value = 0.005058141788724901
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.720897
$ws.Range("H5").Value = 2.162691
$ws.Range("I5").Value = 0.0284720950782092
$ws.Range("J5").Value = 0.02847209507820921
$ws.Range("M5").Value = 0.2545846666666667
$ws.Range("N5").Value = 0.763754
$ws.Range("O5").Value = 0.01165670745182886
$ws.Range("P5").Value = 0.01165670745182886
$ws.Range("Q5").Value = 0.183529322446
$ws.Range("R5").Value = 1.651763902014
$ws.Range("S5").Value = 0.0003318908828673411
$ws.Range("T5").Value = 0.0003318908828673411
$ws.Range("I6").Value = 0.9467537483444692
$ws.Range("J6").Value = 0.9467537483444693
$ws.Range("M6").Value = 0.8908616666666668
$ws.Range("N6").Value = 2.672585
$ws.Range("O6").Value = 0.04079002072021364
$ws.Range("P6").Value = 0.04079002072021363
$ws.Range("Q6").Value = 21.35507485397666
$ws.Range("R6").Value = 192.19567368579
$ws.Range("S6").Value = 0.03861810501191083
$ws.Range("T6").Value = 0.03861810501191083
$ws.Range("I7").Value = 0.9467537483444692
$ws.Range("J7").Value = 0.9467537483444693
$ws.Range("M7").Value = 16.81477433333333
$ws.Range("O7").Value = 0.7699006693471485
$ws.Range("P7").Value = 0.7699006693471484
$ws.Range("Q7").Value = 403.0712937561112
$ws.Range("R7").Value = 3627.641643805001
$ws.Range("S7").Value = 0.7289063445573286
$ws.Range("T7").Value = 0.7289063445573286
$ws.Range("I8").Value = 0.9467537483444692
$ws.Range("J8").Value = 0.9467537483444693
$ws.Range("M8").Value = 3.879966
$ws.Range("N8").Value = 11.639898
$ws.Range("O8").Value = 0.1776526024808091
$ws.Range("P8").Value = 0.1776526024808091
$ws.Range("Q8").Value = 93.00766601722798
$ws.Range("R8").Value = 837.0689941550519
$ws.Range("S8").Value = 0.168193267301856
$ws.Range("T8").Value = 0.168193267301856
$ws.Range("I9").Value = 0.9467537483444692
$ws.Range("J9").Value = 0.9467537483444693
$ws.Range("M9").Value = 0.2545846666666667
$ws.Range("N9").Value = 0.763754
$ws.Range("O9").Value = 0.01165670745182886
$ws.Range("P9").Value = 0.01165670745182886
$ws.Range("Q9").Value = 6.102714727510666
$ws.Range("R9").Value = 54.924432547596
$ws.Range("S9").Value = 0.01103603147337388
$ws.Range("T9").Value = 0.01103603147337388
$ws.Range("G10").Value = 0.6272673333333333
$ws.Range("H10").Value = 1.881802
$ws.Range("I10").Value = 0.0247741565773216
$ws.Range("J10").Value = 0.0247741565773216
$ws.Range("M10").Value = 0.8908616666666668
$ws.Range("N10").Value = 2.672585
$ws.Range("O10").Value = 0.04079002072021364
$ws.Range("P10").Value = 0.04079002072021363
$ws.Range("Q10").Value = 0.5588084220188889
$ws.Range("R10").Value = 5.02927579817
$ws.Range("S10").Value = 0.001010538360114765
$ws.Range("T10").Value = 0.001010538360114765
$ws.Range("G11").Value = 0.6272673333333333
$ws.Range("H11").Value = 1.881802
$ws.Range("I11").Value = 0.0247741565773216
$ws.Range("J11").Value = 0.0247741565773216
$ws.Range("M11").Value = 16.81477433333333
$ws.Range("O11").Value = 0.7699006693471485
$ws.Range("P11").Value = 0.7699006693471484
$ws.Range("Q11").Value = 10.54735865667178
$ws.Range("R11").Value = 94.92622791004599
$ws.Range("S11").Value = 0.01907363973139096
$ws.Range("T11").Value = 0.01907363973139096
$ws.Range("G12").Value = 0.6272673333333333
$ws.Range("H12").Value = 1.881802
$ws.Range("I12").Value = 0.0247741565773216
$ws.Range("J12").Value = 0.0247741565773216
$ws.Range("M12").Value = 3.879966
$ws.Range("N12").Value = 11.639898
$ws.Range("O12").Value = 0.1776526024808091
$ws.Range("P12").Value = 0.1776526024808091
$ws.Range("Q12").Value = 2.433775926244
$ws.Range("R12").Value = 21.903983336196
$ws.Range("S12").Value = 0.004401193390228236
$ws.Range("T12").Value = 0.004401193390228236
$ws.Range("G13").Value = 0.6272673333333333
$ws.Range("H13").Value = 1.881802
$ws.Range("I13").Value = 0.0247741565773216
$ws.Range("J13").Value = 0.0247741565773216
$ws.Range("M13").Value = 0.2545846666666667
$ws.Range("N13").Value = 0.763754
$ws.Range("O13").Value = 0.01165670745182886
$ws.Range("P13").Value = 0.01165670745182886
$ws.Range("Q13").Value = 0.1596926449675556
$ws.Range("R13").Value = 1.437233804708
$ws.Range("S13").Value = 0.0002887850955876397
$ws.Range("T13").Value = 0.0002887850955876397
